# New crime data collected - weekly CompStat update for the 115th Precinct.
# Updates the report header (volume/week-number + reporting date range) and
# refreshes every statistic in the crime-complaints table (rows 14-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 32   Number  48" -> "...Number  49" --------------------
$c = $ws.Range("A8")
$full = $c.Text
$idx = $full.LastIndexOf("48")
$c.Characters($idx + 1, 2).Text = "49"

# --- Header: report week "11/24/2025 Through 11/30/2025" -> "12/1/2025 Through 12/7/2025"
$c = $ws.Range("C9")
$full = $c.Text
$idx = $full.IndexOf("11/24/2025")
$c.Characters($idx + 1, 10).Text = "12/1/2025"
$full = $c.Text
$idx = $full.IndexOf("11/30/2025")
$c.Characters($idx + 1, 10).Text = "12/7/2025"

# --- Cells that flip between a numeric value and the "N/A" text markers -----
# ("0" = shared string used when a count collapses to n/a, "***.*" = shared
# string used when a percentage is undefined). We copy an existing cell that
# already carries the desired style+content so the destination picks up the
# exact same style index / shared-string reference, then (when the target is
# numeric) overwrite just the value.

# -> becomes "N/A" (style 13, text "0", shared string used by D14)
$ws.Range("D14").Copy($ws.Range("C15"))
$ws.Range("D14").Copy($ws.Range("C22"))
$ws.Range("D14").Copy($ws.Range("C27"))
$ws.Range("D14").Copy($ws.Range("D29"))
$ws.Range("D14").Copy($ws.Range("D30"))

# -> becomes "N/A" (style 13, text "***.*", shared string used by E14)
$ws.Range("E14").Copy($ws.Range("E29"))
$ws.Range("E14").Copy($ws.Range("E30"))

# -> becomes a plain number (style 14), donor F22 keeps its own value untouched
$ws.Range("F22").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("F22").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 3

# -> becomes a plain number (style 15), donor H22 keeps its own value untouched
$ws.Range("H22").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("H22").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -66.666666666666

# --- Remaining numeric-only refreshes across the table ----------------------
$ws.Range("L14").Value = 50
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 400
$ws.Range("M15").Value = 3.448275862068
$ws.Range("N15").Value = -9.090909090909
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 33
$ws.Range("H16").Value = -72.727272727272
$ws.Range("I16").Value = 283
$ws.Range("J16").Value = 409
$ws.Range("K16").Value = -30.806845965770
$ws.Range("L16").Value = -15.269461077844
$ws.Range("M16").Value = -12.654320987654
$ws.Range("N16").Value = -76.860179885527
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 53
$ws.Range("H17").Value = -50.943396226415
$ws.Range("I17").Value = 510
$ws.Range("J17").Value = 613
$ws.Range("K17").Value = -16.802610114192
$ws.Range("L17").Value = 14.093959731543
$ws.Range("M17").Value = 65.584415584415
$ws.Range("N17").Value = 25.615763546798
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -28.571428571428
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 145
$ws.Range("J18").Value = 207
$ws.Range("K18").Value = -29.951690821256
$ws.Range("L18").Value = 9.848484848484
$ws.Range("M18").Value = -47.841726618705
$ws.Range("N18").Value = -92.195909580193
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -42.105263157894
$ws.Range("G19").Value = 74
$ws.Range("H19").Value = -29.729729729729
$ws.Range("I19").Value = 683
$ws.Range("J19").Value = 910
$ws.Range("K19").Value = -24.945054945054
$ws.Range("L19").Value = -14.943960149439
$ws.Range("M19").Value = 44.092827004219
$ws.Range("N19").Value = -49.594095940959
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 200
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = -22.727272727272
$ws.Range("I20").Value = 230
$ws.Range("J20").Value = 291
$ws.Range("K20").Value = -20.962199312714
$ws.Range("L20").Value = -29.012345679012
$ws.Range("M20").Value = 0.436681222707
$ws.Range("N20").Value = -89.115002366303
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 48
$ws.Range("E21").Value = -45.833333333333
$ws.Range("F21").Value = 117
$ws.Range("G21").Value = 199
$ws.Range("H21").Value = -41.206030150753
$ws.Range("I21").Value = 1887
$ws.Range("J21").Value = 2474
$ws.Range("K21").Value = -23.726758286176
$ws.Range("L21").Value = -8.708272859216
$ws.Range("M21").Value = 14.711246200607
$ws.Range("N21").Value = -73.062098501070
$ws.Range("J22").Value = 75
$ws.Range("K22").Value = -44
$ws.Range("L22").Value = -49.397590361445
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 8.333333333333
$ws.Range("F24").Value = 84
$ws.Range("G24").Value = 120
$ws.Range("H24").Value = -30
$ws.Range("I24").Value = 1312
$ws.Range("J24").Value = 1864
$ws.Range("K24").Value = -29.613733905579
$ws.Range("L24").Value = -30.175625332623
$ws.Range("M24").Value = 22.846441947565
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 23
$ws.Range("H25").Value = -60.344827586206
$ws.Range("I25").Value = 428
$ws.Range("J25").Value = 954
$ws.Range("K25").Value = -55.136268343815
$ws.Range("L25").Value = -55.042016806722
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 29
$ws.Range("E26").Value = -41.379310344827
$ws.Range("F26").Value = 72
$ws.Range("G26").Value = 94
$ws.Range("H26").Value = -23.404255319148
$ws.Range("I26").Value = 976
$ws.Range("J26").Value = 1243
$ws.Range("K26").Value = -21.480289621882
$ws.Range("L26").Value = 5.742145178764
$ws.Range("M26").Value = 12.313003452244
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 150
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 150
$ws.Range("I28").Value = 113
$ws.Range("J28").Value = 134
$ws.Range("K28").Value = -15.671641791044
$ws.Range("L28").Value = -18.115942028985
$ws.Range("N29").Value = -98.245614035087
$ws.Range("N30").Value = -98.113207547169
